$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155, shifting existing rows 155-182 down to 156-183.
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row 155 with the new data record.
$ws.Range('A155').Value = 5
$ws.Range('B155').Value = 'Macroferia Regional de Talca'
$ws.Range('C155').Value = 'Maule'
$ws.Range('D155').Value = 44491
$ws.Range('E155').Value = 7
$ws.Range('F155').Value = 100112006
$ws.Range('G155').Value = 'Repollo'
$ws.Range('H155').Value = 'Crespo record'
$ws.Range('I155').Value = 'Primera'
$ws.Range('J155').Value = 5000
$ws.Range('K155').Value = 700
$ws.Range('L155').Value = 700
$ws.Range('M155').Value = 700
$ws.Range('N155').Value = '$/unidad'
$ws.Range('O155').Value = 'Provincia del Elquí'
$ws.Range('P155').Value = 700
$ws.Range('Q155').Value = 1
$ws.Range('R155').Value = 'Hortaliza'
